$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4: column I (praclen) goes from 4 to 5
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 5

# The old row 5 data (A5..J5) is being pushed down to become the new row 6,
# with its I value also updated from 4 to 5. Do this first before overwriting row 5.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 31
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim1_2"

# New row 5 data replaces the old row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 61
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = "train_dim1_2"
